$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at row 11, pushing existing rows 11-15 down to 12-16.
$ws.Rows.Item(11).Insert()

# Copy style of row 12 (old row 11, which had the date number format) to new row 11's D cell.
$ws.Cells.Item(12, 4).Copy()
$ws.Cells.Item(11, 4).PasteSpecial(-4122) | Out-Null

# Fill the new row 11 with the weekly record.
$ws.Cells.Item(11, 1).Value = 8
$ws.Cells.Item(11, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(11, 3).Value = "Coquimbo"
$ws.Cells.Item(11, 4).Value = 44985
$ws.Cells.Item(11, 5).Value = 4
$ws.Cells.Item(11, 6).Value = 100112039
$ws.Cells.Item(11, 7).Value = "Ciboulette"
$ws.Cells.Item(11, 8).Value = "Sin especificar"
$ws.Cells.Item(11, 9).Value = "Primera"
$ws.Cells.Item(11, 10).Value = 1000
$ws.Cells.Item(11, 11).Value = 2000
$ws.Cells.Item(11, 12).Value = 2500
$ws.Cells.Item(11, 13).Value = 2250
$ws.Cells.Item(11, 14).Value = "`$/docena de atados"
$ws.Cells.Item(11, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(11, 16).Value = 750
$ws.Cells.Item(11, 17).Value = 3
$ws.Cells.Item(11, 18).Value = "Hortaliza"
